$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2025-08-02 Saturday"

# Update the division problems in the table, addressed by row/column
# so duplicate expressions (e.g. "95÷9=") are each replaced independently
# instead of relying on Find, which always matches from the start of the document.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "64÷6="
$tbl.Cell(1, 2).Range.Text = "48÷4="
$tbl.Cell(1, 3).Range.Text = "66÷6="
$tbl.Cell(1, 4).Range.Text = "52÷9="
$tbl.Cell(1, 5).Range.Text = "41÷7="

$tbl.Cell(5, 1).Range.Text = "28÷3="
$tbl.Cell(5, 2).Range.Text = "37÷2="
$tbl.Cell(5, 3).Range.Text = "37÷5="
$tbl.Cell(5, 4).Range.Text = "29÷9="
$tbl.Cell(5, 5).Range.Text = "36÷8="

$tbl.Cell(9, 1).Range.Text = "77÷5="
$tbl.Cell(9, 2).Range.Text = "78÷3="
$tbl.Cell(9, 3).Range.Text = "85÷4="
$tbl.Cell(9, 4).Range.Text = "65÷5="
$tbl.Cell(9, 5).Range.Text = "35÷6="

$tbl.Cell(13, 1).Range.Text = "56÷4="
$tbl.Cell(13, 2).Range.Text = "91÷5="
$tbl.Cell(13, 3).Range.Text = "65÷8="
$tbl.Cell(13, 4).Range.Text = "49÷7="
$tbl.Cell(13, 5).Range.Text = "54÷4="

$tbl.Cell(17, 1).Range.Text = "63÷6="
$tbl.Cell(17, 2).Range.Text = "37÷2="
$tbl.Cell(17, 3).Range.Text = "31÷3="
$tbl.Cell(17, 4).Range.Text = "38÷3="
$tbl.Cell(17, 5).Range.Text = "82÷8="
